$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.598387956619263
$ws.Range("B1").Value = 2.387347936630249
$ws.Range("C1").Value = 2.71203088760376
$ws.Range("D1").Value = 3.152851104736328
$ws.Range("E1").Value = 1.934929013252258
